$wb = $excel.ActiveWorkbook

# Sheet "建物" (Building): property_category column (I) changes from "land" to "building" for data rows
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"

# Sheet "汽車" (Car): property_category column (H) changes from "land" to "car" for data rows
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
$wsCar.Range("H3").Value = "car"
